# IST price update 2025-12-22 19:03
#
# A new price check was run. This inserts a brand new "latest" timestamp
# column at column B (shifting every existing timestamp/price column one
# place to the right, B->C, C->D, ... AB->AC) and records the current
# price snapshot in the new column B. Where a SKU has no explicit price
# recorded in its (old) first data column, the most recent previously
# known price is carried forward into the new column, exactly like the
# rest of the sheet already does for missing checks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 26
$lastColBeforeInsert = 28   # column AB

# Capture the "latest known price" for every SKU row *before* inserting,
# i.e. the first numeric value found scanning from column B rightwards.
# This is what should be carried into the freshly inserted column.
$carryForward = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $val = $null
    for ($c = 2; $c -le $lastColBeforeInsert; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -ne $null -and $cell.Value2 -ne "") {
            $val = $cell.Value2
            break
        }
    }
    $carryForward[$r] = $val
}

# Insert a new column before column B; this shifts columns B:AB to C:AC.
$ws.Columns("B:B").Insert()

# The freshly inserted column loses its explicit width definition -
# restore it to match its neighbours (all data columns are width 21).
$ws.Columns("B:B").ColumnWidth = $ws.Columns("C:C").ColumnWidth

# New timestamp header for this price check.
$ws.Range("B1").Value = "2025-12-23 00:28"

# Fill in today's price snapshot (carried forward where unknown).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $carryForward[$r]
}
